$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (Leve Item ID 27772)
$ws.Cells.Item(28, 8).Value = 528.46155
$ws.Cells.Item(28, 9).Value = 489.16666
$ws.Cells.Item(28, 10).Value = 1000
$ws.Cells.Item(28, 11).Value = 489.16666
$ws.Cells.Item(28, 12).Value = 1000
$ws.Cells.Item(28, 13).Value = -4.166659999999979
$ws.Cells.Item(28, 14).Value = -1970

# Row 33 (Leve Item ID 5512)
$ws.Cells.Item(33, 8).Value = 690.3200000000001
$ws.Cells.Item(33, 9).Value = 680.7692
$ws.Cells.Item(33, 11).Value = 680.7692
$ws.Cells.Item(33, 13).Value = -451.7692

# Row 87 (Leve Item ID 10651)
$ws.Cells.Item(87, 8).Value = 26464.76
$ws.Cells.Item(87, 10).Value = 26464.76
$ws.Cells.Item(87, 12).Value = 26464.76
$ws.Cells.Item(87, 14).Value = -28960.76

# Row 90 (Leve Item ID 10651)
$ws.Cells.Item(90, 8).Value = 26464.76
$ws.Cells.Item(90, 10).Value = 26464.76
$ws.Cells.Item(90, 12).Value = 79394.28
$ws.Cells.Item(90, 14).Value = -91874.28

# Row 123 (Leve Item ID 34090)
$ws.Cells.Item(123, 8).Value = 34557.145
$ws.Cells.Item(123, 10).Value = 34557.145
$ws.Cells.Item(123, 12).Value = 34557.145
$ws.Cells.Item(123, 14).Value = -44357.145

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Cells.Item(32, 8).Value = 8321.456
$ws.Cells.Item(32, 9).Value = 6788.406
$ws.Cells.Item(32, 10).Value = 18899.5
$ws.Cells.Item(32, 11).Value = 6788.406
$ws.Cells.Item(32, 12).Value = 18899.5
$ws.Cells.Item(32, 13).Value = -6501.406
$ws.Cells.Item(32, 14).Value = -19473.5

# Row 61 (Leve Item ID 43999)
$ws.Cells.Item(61, 8).Value = 7076.75
$ws.Cells.Item(61, 9).Value = 12170.777
$ws.Cells.Item(61, 10).Value = 2908.9092
$ws.Cells.Item(61, 11).Value = 12170.777
$ws.Cells.Item(61, 12).Value = 2908.9092
$ws.Cells.Item(61, 13).Value = -11958.777
$ws.Cells.Item(61, 14).Value = -3332.9092

# Row 74 (Leve Item ID 44000)
$ws.Cells.Item(74, 8).Value = 1780.7858
$ws.Cells.Item(74, 9).Value = 1560.9166
$ws.Cells.Item(74, 10).Value = 3100
$ws.Cells.Item(74, 11).Value = 1560.9166
$ws.Cells.Item(74, 12).Value = 3100
$ws.Cells.Item(74, 13).Value = -686.9166
$ws.Cells.Item(74, 14).Value = -4848

# Row 77 (Leve Item ID 44000)
$ws.Cells.Item(77, 8).Value = 1780.7858
$ws.Cells.Item(77, 9).Value = 1560.9166
$ws.Cells.Item(77, 10).Value = 3100
$ws.Cells.Item(77, 11).Value = 7804.583000000001
$ws.Cells.Item(77, 12).Value = 15500
$ws.Cells.Item(77, 13).Value = -3436.583000000001
$ws.Cells.Item(77, 14).Value = -24236

# Row 113 (Leve Item ID 26002)
$ws.Cells.Item(113, 8).Value = 31800
$ws.Cells.Item(113, 10).Value = 31800
$ws.Cells.Item(113, 12).Value = 31800
$ws.Cells.Item(113, 14).Value = -40478

# Row 132 (Leve Item ID 43997)
$ws.Cells.Item(132, 8).Value = 3276.575
$ws.Cells.Item(132, 9).Value = 3320.3076
$ws.Cells.Item(132, 10).Value = 3195.3572
$ws.Cells.Item(132, 11).Value = 9960.9228
$ws.Cells.Item(132, 12).Value = 9586.071599999999
$ws.Cells.Item(132, 13).Value = -7430.9228
$ws.Cells.Item(132, 14).Value = -14646.0716

# Row 136 (Leve Item ID 43999)
$ws.Cells.Item(136, 8).Value = 7076.75
$ws.Cells.Item(136, 9).Value = 12170.777
$ws.Cells.Item(136, 10).Value = 2908.9092
$ws.Cells.Item(136, 11).Value = 36512.331
$ws.Cells.Item(136, 12).Value = 8726.7276
$ws.Cells.Item(136, 13).Value = -33962.331
$ws.Cells.Item(136, 14).Value = -13826.7276

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (Leve Item ID 19943)
$ws.Cells.Item(99, 8).Value = 2908.3333
$ws.Cells.Item(99, 9).Value = 2275.7144
$ws.Cells.Item(99, 10).Value = 3310.9092
$ws.Cells.Item(99, 11).Value = 2275.7144
$ws.Cells.Item(99, 12).Value = 3310.9092
$ws.Cells.Item(99, 13).Value = -777.7143999999998
$ws.Cells.Item(99, 14).Value = -6306.9092

# Row 134 (Leve Item ID 43998)
$ws.Cells.Item(134, 8).Value = 4921.032
$ws.Cells.Item(134, 9).Value = 4236.96
$ws.Cells.Item(134, 10).Value = 7771.3335
$ws.Cells.Item(134, 11).Value = 12710.88
$ws.Cells.Item(134, 12).Value = 23314.0005
$ws.Cells.Item(134, 13).Value = -10175.88
$ws.Cells.Item(134, 14).Value = -28384.0005

$ws = $wb.Worksheets.Item("CRP")
# Row 58 (Leve Item ID 44021)
$ws.Cells.Item(58, 8).Value = 11114832
$ws.Cells.Item(58, 9).Value = 1845.8518
$ws.Cells.Item(58, 10).Value = 27784310
$ws.Cells.Item(58, 11).Value = 1845.8518
$ws.Cells.Item(58, 12).Value = 27784310
$ws.Cells.Item(58, 13).Value = -1642.8518
$ws.Cells.Item(58, 14).Value = -27784716

# Row 132 (Leve Item ID 44019)
$ws.Cells.Item(132, 8).Value = 2350.6287
$ws.Cells.Item(132, 9).Value = 2035.96
$ws.Cells.Item(132, 10).Value = 3137.3
$ws.Cells.Item(132, 11).Value = 6107.88
$ws.Cells.Item(132, 12).Value = 9411.900000000001
$ws.Cells.Item(132, 13).Value = -3577.88
$ws.Cells.Item(132, 14).Value = -14471.9

# Row 134 (Leve Item ID 44020)
$ws.Cells.Item(134, 8).Value = 20099.5
$ws.Cells.Item(134, 9).Value = 25661.812
$ws.Cells.Item(134, 10).Value = 5266.6665
$ws.Cells.Item(134, 11).Value = 76985.436
$ws.Cells.Item(134, 12).Value = 15799.9995
$ws.Cells.Item(134, 13).Value = -74450.436
$ws.Cells.Item(134, 14).Value = -20869.9995

# Row 136 (Leve Item ID 44021)
$ws.Cells.Item(136, 8).Value = 11114832
$ws.Cells.Item(136, 9).Value = 1845.8518
$ws.Cells.Item(136, 10).Value = 27784310
$ws.Cells.Item(136, 11).Value = 5537.555399999999
$ws.Cells.Item(136, 12).Value = 83352930
$ws.Cells.Item(136, 13).Value = -2987.555399999999
$ws.Cells.Item(136, 14).Value = -83358030

$ws = $wb.Worksheets.Item("CUL")
# Row 119 (Leve Item ID 27873)
$ws.Cells.Item(119, 8).Value = 3539.4443
$ws.Cells.Item(119, 9).Value = 1997.5
$ws.Cells.Item(119, 10).Value = 3980
$ws.Cells.Item(119, 11).Value = 5992.5
$ws.Cells.Item(119, 12).Value = 11940
$ws.Cells.Item(119, 13).Value = -1154.5
$ws.Cells.Item(119, 14).Value = -21616

$ws = $wb.Worksheets.Item("GSM")
# Row 26 (Leve Item ID 4254)
$ws.Cells.Item(26, 8).Value = 80042
$ws.Cells.Item(26, 10).Value = 80042
$ws.Cells.Item(26, 12).Value = 80042
$ws.Cells.Item(26, 14).Value = -80602

# Row 50 (Leve Item ID 4254)
$ws.Cells.Item(50, 8).Value = 80042
$ws.Cells.Item(50, 10).Value = 80042
$ws.Cells.Item(50, 12).Value = 80042
$ws.Cells.Item(50, 14).Value = -81038

# Row 97 (Leve Item ID 19940)
$ws.Cells.Item(97, 8).Value = 1233.875
$ws.Cells.Item(97, 9).Value = 721.1111
$ws.Cells.Item(97, 10).Value = 1893.1428
$ws.Cells.Item(97, 11).Value = 721.1111
$ws.Cells.Item(97, 12).Value = 1893.1428
$ws.Cells.Item(97, 13).Value = -225.1111
$ws.Cells.Item(97, 14).Value = -2885.1428

# Row 122 (Leve Item ID 36182)
$ws.Cells.Item(122, 8).Value = 7445.3335
$ws.Cells.Item(122, 9).Value = 6000
$ws.Cells.Item(122, 11).Value = 18000
$ws.Cells.Item(122, 13).Value = -15550

# Row 123 (Leve Item ID 34150)
$ws.Cells.Item(123, 8).Value = 21131.285
$ws.Cells.Item(123, 10).Value = 21131.285
$ws.Cells.Item(123, 12).Value = 21131.285
$ws.Cells.Item(123, 14).Value = -26031.285

# Row 132 (Leve Item ID 44008)
$ws.Cells.Item(132, 8).Value = 2678.638
$ws.Cells.Item(132, 9).Value = 2304.4324
$ws.Cells.Item(132, 11).Value = 6913.297200000001
$ws.Cells.Item(132, 13).Value = -4383.297200000001

# Row 136 (Leve Item ID 42218)
$ws.Cells.Item(136, 8).Value = 6362.636
$ws.Cells.Item(136, 10).Value = 6362.636
$ws.Cells.Item(136, 12).Value = 19087.908
$ws.Cells.Item(136, 14).Value = -24187.908

$ws = $wb.Worksheets.Item("LTW")
# Row 122 (Leve Item ID 36247)
$ws.Cells.Item(122, 8).Value = 3692.6155
$ws.Cells.Item(122, 9).Value = 2729.1428
$ws.Cells.Item(122, 10).Value = 4816.6665
$ws.Cells.Item(122, 11).Value = 8187.428400000001
$ws.Cells.Item(122, 12).Value = 14449.9995
$ws.Cells.Item(122, 13).Value = -5737.428400000001
$ws.Cells.Item(122, 14).Value = -19349.9995

# Row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, 8).Value = 3593.6086
$ws.Cells.Item(132, 9).Value = 3274.6667
$ws.Cells.Item(132, 11).Value = 9824.000100000001
$ws.Cells.Item(132, 13).Value = -7294.000100000001

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (Leve Item ID 27746)
$ws.Cells.Item(107, 8).Value = 1855
$ws.Cells.Item(107, 9).Value = 525.5
$ws.Cells.Item(107, 11).Value = 1576.5
$ws.Cells.Item(107, 13).Value = 343.5

# Row 122 (Leve Item ID 36208)
$ws.Cells.Item(122, 8).Value = 419465.4
$ws.Cells.Item(122, 9).Value = 716500.8
$ws.Cells.Item(122, 10).Value = 3615.9
$ws.Cells.Item(122, 11).Value = 2149502.4
$ws.Cells.Item(122, 12).Value = 10847.7
$ws.Cells.Item(122, 13).Value = -2147052.4
$ws.Cells.Item(122, 14).Value = -15747.7

# Row 136 (Leve Item ID 44031)
$ws.Cells.Item(136, 8).Value = 1612.0344
$ws.Cells.Item(136, 9).Value = 867.58826
$ws.Cells.Item(136, 10).Value = 2666.6667
$ws.Cells.Item(136, 11).Value = 2602.76478
$ws.Cells.Item(136, 12).Value = 8000.000100000001
$ws.Cells.Item(136, 13).Value = -52.76477999999997
$ws.Cells.Item(136, 14).Value = -13100.0001
